{"js": "// Apply the Astro Babes copy edits using the Word JavaScript API.\nconst replacements = [\n  {\n    from: \"Play Astro Babes Online Slot Game for Free - Review\",\n    to: \"Play Astro Babes Free - Exciting Space-Themed Slot Game\",\n  },\n  {\n    from: \"Attractive graphics and game atmosphere\",\n    to: \"Attractive space-themed graphic design\",\n  },\n  {\n    from: \"Five fascinating protagonists offering winning opportunities\",\n    to: \"Fascinating protagonists with unique values\",\n  },\n  {\n    from: \"High payout potential\",\n    to: \"High potential for significant wins\",\n  },\n  {\n    from: \"Medium-high variance may require patience for significant wins\",\n    to: \"Medium-high variance may require patience\",\n  },\n  {\n    from: \"No progressive jackpot\",\n    to: \"Some players may find the game's theme unappealing\",\n  },\n  {\n    from: \"Read our review of Astro Babes online slot game, play for free, and discover exciting bonus functions, free spins, and high payout potential.\",\n    to: \"Read our review of Astro Babes, a space-themed slot game with exciting bonus functions. Play for free and win big!\",\n  },\n];\n\nfor (const { from, to } of replacements) {\n  const results = context.document.body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the Astro Babes copy edits using the Word COM object model.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Play Astro Babes Online Slot Game for Free - Review\"; Replace = \"Play Astro Babes Free - Exciting Space-Themed Slot Game\" },\n    @{ Find = \"Attractive graphics and game atmosphere\"; Replace = \"Attractive space-themed graphic design\" },\n    @{ Find = \"Five fascinating protagonists offering winning opportunities\"; Replace = \"Fascinating protagonists with unique values\" },\n    @{ Find = \"High payout potential\"; Replace = \"High potential for significant wins\" },\n    @{ Find = \"Medium-high variance may require patience for significant wins\"; Replace = \"Medium-high variance may require patience\" },\n    @{ Find = \"No progressive jackpot\"; Replace = \"Some players may find the game's theme unappealing\" },\n    @{ Find = \"Read our review of Astro Babes online slot game, play for free, and discover exciting bonus functions, free spins, and high payout potential.\"; Replace = \"Read our review of Astro Babes, a space-themed slot game with exciting bonus functions. Play for free and win big!\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    while ($find.Execute()) {\n        # Assign Range.Text directly (rather than Find.Replacement + Execute's\n        # Replace argument) so AutoCorrect's smart-quote substitution never\n        # touches the straight apostrophe in the replacement copy.\n        $range.Text = $r.Replace\n        $range.Collapse(0)\n        $find = $range.Find\n        $find.ClearFormatting()\n        $find.Text = $r.Find\n        $find.Forward = $true\n        $find.Wrap = 1\n        $find.MatchCase = $true\n        $find.MatchWholeWord = $false\n        $find.MatchWildcards = $false\n    }\n}\n"}
